# Update docx templates for correct bibliography paragraph.
#
# 1. The lone paragraph in the template currently uses the "Caption"
#    style; it should use the "Bibliography" style instead.
# 2. The "Bibliography" paragraph style itself needs a hanging-indent
#    paragraph format (left indent 1134 twips / 56.7pt, hanging 1134
#    twips / 56.7pt -> FirstLineIndent = -56.7pt).
# 3. The "Footnote Text" paragraph style needs to be justified
#    (both-sides alignment).

$d = $word.ActiveDocument

# 1) Re-style the template's single paragraph: Caption -> Bibliography.
$d.Paragraphs(1).Style = "Bibliography"

# 2) Give the Bibliography style a hanging indent of 1134 twips
#    (1134 / 20 = 56.7 points), matching <w:ind w:left="1134" w:hanging="1134"/>.
$biblioStyle = $d.Styles("Bibliography")
$biblioStyle.ParagraphFormat.LeftIndent = 56.7
$biblioStyle.ParagraphFormat.FirstLineIndent = -56.7

# 3) Justify the Footnote Text style (both-sides alignment -> <w:jc w:val="both"/>).
$footnoteStyle = $d.Styles("Footnote Text")
$footnoteStyle.ParagraphFormat.Alignment = 3
